$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'49.867.26"
$ws.Range("E2").Value = "  +3.56%  "

$ws.Range("D3").Value = "'2.625.22"
$ws.Range("E3").Value = "  +4.98%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.21%  "

$ws.Range("B5").Value = "Solana"
$ws.Range("C5").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D5").Value = "'110.39"
$ws.Range("E5").Value = "  +2.11%  "

$ws.Range("B6").Value = "BNB"
$ws.Range("C6").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D6").Value = "'323.95"
$ws.Range("E6").Value = "  +0.87%  "

$ws.Range("E7").Value = "  +1.89%  "

$ws.Range("E8").Value = "  -0.09%  "

$ws.Range("E9").Value = "  +4.52%  "

$ws.Range("D10").Value = "'41.02"
$ws.Range("E10").Value = "  +2.97%  "

$ws.Range("D11").Value = "'20.66"
$ws.Range("E11").Value = "  +2.74%  "

$ws.Range("E12").Value = "  +1.40%  "

$ws.Range("E13").Value = "  +0.70%  "

$ws.Range("E14").Value = "  +2.70%  "

$ws.Range("D15").Value = "'3.036.96"
$ws.Range("E15").Value = "  +4.86%  "

$ws.Range("D16").Value = "'2.604.29"
$ws.Range("E16").Value = "  +3.97%  "

$ws.Range("D17").Value = "'0.875"
$ws.Range("E17").Value = "  +3.85%  "

$ws.Range("D18").Value = "'49.833.23"
$ws.Range("E18").Value = "  +3.81%  "

$ws.Range("D19").Value = "'3.11"
$ws.Range("E19").Value = "  +12.56%  "

$ws.Range("E20").Value = "  +2.49%  "

$ws.Range("D21").Value = "'6.80"
$ws.Range("E21").Value = "  +0.67%  "

$ws.Range("D22").Value = "'0.0₃0957"
$ws.Range("E22").Value = "  +1.32%  "

$ws.Range("D23").Value = "'282.29"
$ws.Range("E23").Value = "  +1.86%  "

$ws.Range("D24").Value = "'72.92"
$ws.Range("E24").Value = "  +1.37%  "

$ws.Range("D25").Value = "'2.58"
$ws.Range("E25").Value = "  +0.94%  "

$ws.Range("D26").Value = "'26.77"
$ws.Range("E26").Value = "  +4.03%  "

$ws.Range("E27").Value = "  +0.03%  "

$ws.Range("E28").Value = "  -0.90%  "

$ws.Range("E29").Value = "  +1.92%  "

$ws.Range("E30").Value = "  +3.84%  "

$ws.Range("D31").Value = "'36.30"
$ws.Range("E31").Value = "  +2.77%  "

$ws.Range("D32").Value = "'49.62"
$ws.Range("E32").Value = "  +1.16%  "

$ws.Range("D33").Value = "'19.74"
$ws.Range("E33").Value = "  +1.33%  "

$ws.Range("E34").Value = "  +2.79%  "

$ws.Range("E35").Value = "  -0.21%  "

$ws.Range("D36").Value = "'0.0797"
$ws.Range("E36").Value = "  +1.78%  "

$ws.Range("D37").Value = "'2.07"
$ws.Range("E37").Value = "  +5.85%  "

$ws.Range("D38").Value = "'4.77"
$ws.Range("E38").Value = "  +2.94%  "

$ws.Range("D39").Value = "'3.09"
$ws.Range("E39").Value = "  +5.53%  "

$ws.Range("B40").Value = "EnergySwap"
$ws.Range("C40").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D40").Value = "'23.03"
$ws.Range("E40").Value = "  +7.95%  "

$ws.Range("D41").Value = "'123.83"
$ws.Range("E41").Value = "  +2.02%  "

$ws.Range("B42").Value = "Stellar"
$ws.Range("C42").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D42").Value = "'0.113"
$ws.Range("E42").Value = "  +1.20%  "

$ws.Range("E43").Value = "  +0.29%  "

$ws.Range("E44").Value = "  +4.26%  "

$ws.Range("E45").Value = "  +6.31%  "

$ws.Range("D46").Value = "'2.056.19"
$ws.Range("E46").Value = "  +2.73%  "

$ws.Range("D47").Value = "'2.20"
$ws.Range("E47").Value = "  +10.53%  "

$ws.Range("D48").Value = "'2.03"
$ws.Range("E48").Value = "  +10.40%  "

$ws.Range("D49").Value = "'9.09"
$ws.Range("E49").Value = "  +1.47%  "

$ws.Range("E50").Value = "  +4.08%  "

$ws.Range("D51").Value = "'82.23"
$ws.Range("E51").Value = "  +2.96%  "
